# "maf formatting to string"
#
# The MetaboBank MAF NMR template had a handful of "Last updated" /
# "Revision history" dates stored as real Excel date serials (with a
# yyyy-mm-dd number format). This edit converts those cells (and, for
# consistency, every other populated cell on the README sheet plus the
# header cells on MB_MAF) to plain Text-formatted cells, clears out a
# few now-unused blank helper cells/rows, and bumps the "Last updated"
# date to a fresh value.

$wb = $excel.ActiveWorkbook
$readme = $wb.Worksheets.Item("README")
$maf = $wb.Worksheets.Item("MB_MAF")

# ---------------------------------------------------------------------
# README sheet
# ---------------------------------------------------------------------

# A9, A11/B11, B12 and B17 only ever held formatting for blank spacer
# cells - drop them completely so the rows collapse away.
$readme.Range("A9").Clear()
$readme.Range("A11").Clear()
$readme.Range("B11").Clear()
$readme.Range("B12").Clear()
$readme.Range("B17").Clear()

# Every remaining populated cell on the sheet moves to a Text number
# format ("@") instead of General/date.
$readmeCells = @(
    "A1","A2","A3","A4","A5","A6","A7","A8",
    "A10","B10",
    "A12",
    "A13","B13",
    "A14","B14",
    "A15","B15",
    "A16","B16"
)
foreach ($addr in $readmeCells) {
    $readme.Range($addr).NumberFormat = "@"
}

# The date-serial cells become literal yyyy-mm-dd text. B16 kept its
# original "metabolite_category to metabolite_class" date and B13/B14/
# B15 kept theirs; B10 ("Last updated") is bumped forward to the date
# of this edit.
$readme.Range("B16").Value = "2022-12-22"
$readme.Range("B13").Value = "2021-11-09"
$readme.Range("B14").Value = "2022-03-25"
$readme.Range("B15").Value = "2022-07-21"
$readme.Range("B10").Value = "2024-03-12"

# ---------------------------------------------------------------------
# MB_MAF sheet - same Text-format normalization, values untouched.
# ---------------------------------------------------------------------

$maf.Range("A1").NumberFormat = "@"
$maf.Range("A2").NumberFormat = "@"
$maf.Range("E2").NumberFormat = "@"
$maf.Range("Q2").NumberFormat = "@"
$maf.Range("A3:P3").NumberFormat = "@"
